# Generate Report for Handback
# The c24ac456-c25d-494b-9f17-797937471c65.md file has finished handback
# (it now matches en-US), so flip its status from "Ready for handoff" to
# "Handed back: in sync with en-US" across the Overview, zh-cn and de-de
# sheets, and stamp the new handback datetime for each locale.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet: row 3 is the c24ac456... file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet: row 3 is the c24ac456... file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusHandedBack
$wsZhCn.Range("H3").Value = "2016-03-20 17:08:31"

# --- de-de sheet: row 3 is the c24ac456... file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusHandedBack
$wsDeDe.Range("H3").Value = "2016-03-20 17:08:44"
